# pauta.docx -- "Creation of Schedule working"
#
# 1) Resize the "Frame1" textbox shape slightly (wp:extent / a:ext / VML
#    fallback all grow a touch -- this mirrors Word recalculating the
#    textbox's reported extents).
# 2) Merge "às " + "horaAudiencia" into a single run "às hora " and drop
#    the leading space on the following "- Valor Anterior: " run.
# 3) Split "anoProcessso numeroProcesso} " into a superscript leading
#    space, "ano" and " numeroProcesso ".
# 4) Swap the names of two pairs of zero-width bookmarks.
# 5) Regenerate the per-character Fieldmark bookmarks wrapping the
#    "«Ignorar registro se...»" field result.

$d = $word.ActiveDocument

# --- 1. Resize the Frame1 shape -------------------------------------------
$shp = $d.Shapes.Item(1)
$shp.Width = 37.1
$shp.Height = 17.9

# --- 2. "às horaAudiencia" -> "às hora " -----------------------------------
$d.Content.Find.Execute("às horaAudiencia", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "às hora ", 2)
$d.Content.Find.Execute(" - Valor Anterior: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "- Valor Anterior: ", 2)

# --- 3. "anoProcessso numeroProcesso} " -> superscript space + "ano" + " numeroProcesso " --
$d.Content.Find.Execute("anoProcessso numeroProcesso} ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " ano numeroProcesso ", 2)

$findRng = $d.Content
$findRng.Find.Execute(" ano numeroProcesso ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$anoPos = $findRng.Start

$spaceRng = $d.Range($anoPos, $anoPos + 1)
$spaceRng.Font.Superscript = $true

# Force a run break between "ano" and " numeroProcesso " (identical
# formatting either side, so only a structural edit -- e.g. a transient
# bookmark -- will split the run).
$splitPos = $anoPos + 4
$splitRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("__tmp_split_ano", $splitRng)
$d.Bookmarks.Item("__tmp_split_ano").Delete()

# --- 4. Swap paired bookmark names -----------------------------------------
function Swap-Bookmarks($nameA, $nameB) {
    $bmA = $d.Bookmarks.Item($nameA)
    $bmB = $d.Bookmarks.Item($nameB)
    $rngA = $d.Range($bmA.Start, $bmA.End)
    $rngB = $d.Range($bmB.Start, $bmB.End)
    $bmA.Delete()
    $bmB.Delete()
    $d.Bookmarks.Add($nameA, $rngB)
    $d.Bookmarks.Add($nameB, $rngA)
}

Swap-Bookmarks "__Fieldmark__2_781997288" "__Fieldmark__3_781997288"
Swap-Bookmarks "__DdeLink__73_451291707" "__Fieldmark__4_781997288"

# --- 5. Regenerate the Fieldmark bookmarks around «Ignorar registro se...»» --
$fieldRng = $d.Content
$fieldRng.Find.Execute([char]0xAB + "Ignorar registro se..." + [char]0xBB, $true, $false, `
                        $false, $false, $false, $true, 1, $false, "", 0)
$fieldStart = $fieldRng.Start
$fieldEnd = $fieldRng.End

$oldNames = @("__Fieldmark__555_1701348422", "__Fieldmark__416_1701348422", `
              "__Fieldmark__50_1821181616", "__Fieldmark__59_953294636", `
              "__Fieldmark__66_451291707", "__Fieldmark__6_781997288")
foreach ($n in $oldNames) {
    $d.Bookmarks.Item($n).Delete()
}

$newNames = @("__Fieldmark__49_186955932", "__Fieldmark__340_1145423660", `
              "__Fieldmark__50_333644460", "__Fieldmark__555_1701348422", `
              "__Fieldmark__416_1701348422", "__Fieldmark__50_1821181616", `
              "__Fieldmark__59_953294636", "__Fieldmark__66_451291707", `
              "__Fieldmark__6_781997288")

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $start = $fieldStart + $i
    $d.Bookmarks.Add($newNames[$i], $d.Range($start, $fieldEnd))
}
